$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 686 (pushes the existing 2026/12/29.. data down by one row)
$ws.Rows.Item(686).Insert()

# Force column A to stay plain text so "2026/01/21" is not auto-converted
# into a date serial number, then strip the temporary text format so the
# new cell ends up unstyled like its siblings.
$ws.Range("A686").NumberFormat = "@"
$ws.Range("A686").Value = "2026/01/21"
$ws.Range("A686").ClearFormats()

$ws.Range("B686").Value = "水"
$ws.Range("C686").Value = 16
$ws.Range("D686").Value = 201
